# Update the 取得日時 (acquisition timestamp) column on the "ランサーズ" sheet.
# All data rows (2-17) currently hold "2025-10-10 12:34:56" in column A;
# this run appended/refreshed data at 2025-10-10 12:46:22, so every row's
# timestamp is bumped to the new value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-10 12:46:22"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = $ws.UsedRange.Rows.Count
}

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
